$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trip")

$ws.Cells.Item(17,1).Value = "T15"
$ws.Cells.Item(17,2).Value = "Shared Trip ID"
$ws.Cells.Item(17,3).Value = "String"
$ws.Cells.Item(17,4).Value = "Yes"
$ws.Cells.Item(17,5).Value = "Each complete empty-to-empty run of a vehicle should be assigned a unique ID and this ID should be entered into this field. The ID should be a non-case-sensitive string value with any letters represented in their capital forms. Each TNP may use a coding system of its liking to assign the IDs, although the City of Chicago reserves the right to apply further restrictions, as it finds necessary. Therefore the IDs need only be unique within a TNP, not between TNPs.  Every trip record within the empty-to-empty run should list the same Shared Trip ID."

$ws.Range("A17:B17,D17").VerticalAlignment = -4160
$ws.Range("C17:C17,E17:E17").VerticalAlignment = -4160
$ws.Range("C17:C17,E17:E17").WrapText = $true
